$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RESTATUS description (column C, row 11) to add the new
# status code "4 = ถอนคำร้อง" (withdrawn request). The Constraints cell
# (D11, "Default: 0") is left untouched.
$ws.Range("C11").Value = "สถานะของคำร้อง (0 = รอดำเนินการ, 1 = อยู่ระหว่างดำเนินการ, 2 = เสร็จสิ้น, 3 = ปฏิเสธ, 4 = ถอนคำร้อง)"

# The longer description no longer fits the previous auto (best-fit)
# column width, so the column was widened and manually fixed instead
# of relying on best-fit sizing.
$ws.Columns.Item(3).ColumnWidth = 91.43
